$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellAddress, $text)
    $c = $ws.Range($cellAddress)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.628.43"
$ws.Range("E2").Value = "  +0.07%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.651.02"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "596.71"
$ws.Range("E5").Value = "  -0.68%  "

# Row 6 - Solana
Set-TextValue "D6" "156.32"
$ws.Range("E6").Value = "  -0.20%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +3.78%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +3.99%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.396"
$ws.Range("E10").Value = "  -0.29%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.80"
$ws.Range("E11").Value = "  -1.92%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +1.10%  "

# Row 13 - Avalanche
Set-TextValue "D13" "28.77"
$ws.Range("E13").Value = "  -2.12%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +0.88%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.129.18"
$ws.Range("E15").Value = "  -0.53%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "65.481.08"
$ws.Range("E16").Value = "  +0.15%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.643.96"
$ws.Range("E17").Value = "  -1.10%  "

# Row 18 - Chainlink
Set-TextValue "D18" "12.60"
$ws.Range("E18").Value = "  +0.92%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -1.29%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -1.04%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "348.76"
$ws.Range("E21").Value = "  -0.56%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.11%  "

# Row 23 - Litecoin
Set-TextValue "D23" "69.06"
$ws.Range("E23").Value = "  -0.96%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  +2.70%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("E25").Value = "  +0.17%  "

# Row 26 - SuiNetwork
$ws.Range("E26").Value = "  +1.42%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  -0.26%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -1.99%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.05%  "

# Row 30 - Aptos
$ws.Range("E30").Value = "  -2.31%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.14%  "

# Row 32 - Bittensor
Set-TextValue "D32" "528.63"
$ws.Range("E32").Value = "  -2.17%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  -0.05%  "

# Row 34 - RenderToken
$ws.Range("E34").Value = "  -1.33%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  -0.59%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "  -0.28%  "

# Row 37 - EthereumClassic
Set-TextValue "D37" "20.37"
$ws.Range("E37").Value = "  -0.10%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  +0.02%  "

# Row 39 - Monero
Set-TextValue "D39" "156.20"
$ws.Range("E39").Value = "  -1.76%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  -0.92%  "

# Row 41 - USDe
$ws.Range("E41").Value = "  +0.00%  "

# Row 42 - Aave
Set-TextValue "D42" "161.00"
$ws.Range("E42").Value = "  -2.69%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  +0.10%  "

# Row 44 - Hedera
$ws.Range("E44").Value = "  -0.84%  "

# Row 45 - dogwifhat
Set-TextValue "D45" "2.27"
$ws.Range("E45").Value = "  +0.36%  "

# Row 46 - InjectiveProtocol
Set-TextValue "D46" "22.62"
$ws.Range("E46").Value = "  -1.78%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  -1.81%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  -1.80%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -0.59%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +11.83%  "

# Row 51 - EnergySwap
Set-TextValue "D51" "19.71"
$ws.Range("E51").Value = "  -1.07%  "
